# Auto-generated edit script: updates computed profit/price figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets,
# reflecting refreshed market data from the scheduled runner.

$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1984.6
$ws.Range("J2").Value = 1624.25
$ws.Range("L2").Value = 1624.25
$ws.Range("N2").Value = -1850.25
$ws.Range("H5").Value = 291.88235
$ws.Range("I5").Value = 244.5
$ws.Range("J5").Value = 405.6
$ws.Range("K5").Value = 244.5
$ws.Range("L5").Value = 405.6
$ws.Range("M5").Value = -129.5
$ws.Range("N5").Value = -635.6
$ws.Range("H98").Value = 1602.4073
$ws.Range("I98").Value = 1467.0454
$ws.Range("K98").Value = 1467.0454
$ws.Range("M98").Value = 30.95460000000003
$ws.Range("H113").Value = 136214.06
$ws.Range("I113").Value = 3499.889
$ws.Range("J113").Value = 335285.34
$ws.Range("K113").Value = 3499.889
$ws.Range("L113").Value = 335285.34
$ws.Range("M113").Value = -245.8890000000001
$ws.Range("N113").Value = -341793.34
$ws.Range("H115").Value = 14318287
$ws.Range("I115").Value = 14318287
$ws.Range("K115").Value = 42954861
$ws.Range("M115").Value = -42953294
$ws.Range("H122").Value = 1602.4073
$ws.Range("I122").Value = 1467.0454
$ws.Range("K122").Value = 4401.1362
$ws.Range("M122").Value = -1951.1362
$ws.Range("H129").Value = 1732.3334
$ws.Range("I129").Value = 1205.4615
$ws.Range("K129").Value = 3616.3845
$ws.Range("M129").Value = 1383.6155
$ws.Range("H132").Value = 1343
$ws.Range("I132").Value = 1225.381
$ws.Range("K132").Value = 3676.143
$ws.Range("M132").Value = -1146.143
$ws.Range("H138").Value = 2876.8572
$ws.Range("J138").Value = 3791.762
$ws.Range("L138").Value = 11375.286
$ws.Range("N138").Value = -21655.286

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 196.125
$ws.Range("I5").Value = 305.77777
$ws.Range("J5").Value = 55.142857
$ws.Range("K5").Value = 305.77777
$ws.Range("L5").Value = 55.142857
$ws.Range("M5").Value = -193.77777
$ws.Range("N5").Value = -279.142857
$ws.Range("H32").Value = 13829.469
$ws.Range("I32").Value = 14046.421
$ws.Range("K32").Value = 14046.421
$ws.Range("M32").Value = -13759.421
$ws.Range("H43").Value = 39998
$ws.Range("J43").Value = 39998
$ws.Range("L43").Value = 39998
$ws.Range("N43").Value = -40624
$ws.Range("H45").Value = 39426.5
$ws.Range("I45").Value = 43032.89
$ws.Range("K45").Value = 43032.89
$ws.Range("M45").Value = -42655.89
$ws.Range("H61").Value = 4016.9814
$ws.Range("I61").Value = 3171.5
$ws.Range("J61").Value = 10780.833
$ws.Range("K61").Value = 3171.5
$ws.Range("L61").Value = 10780.833
$ws.Range("M61").Value = -2959.5
$ws.Range("N61").Value = -11204.833
$ws.Range("H74").Value = 317976.6
$ws.Range("I74").Value = 361687.53
$ws.Range("J74").Value = 12000
$ws.Range("K74").Value = 361687.53
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = -360813.53
$ws.Range("N74").Value = -13748
$ws.Range("H77").Value = 317976.6
$ws.Range("I77").Value = 361687.53
$ws.Range("J77").Value = 12000
$ws.Range("K77").Value = 1808437.65
$ws.Range("L77").Value = 60000
$ws.Range("M77").Value = -1804069.65
$ws.Range("N77").Value = -68736
$ws.Range("H97").Value = 1159080.5
$ws.Range("I97").Value = 1278714.4
$ws.Range("K97").Value = 1278714.4
$ws.Range("M97").Value = -1278218.4
$ws.Range("H102").Value = 1277.75
$ws.Range("I102").Value = 1277.75
$ws.Range("K102").Value = 1277.75
$ws.Range("M102").Value = 344.25
$ws.Range("H122").Value = 4453.7827
$ws.Range("I122").Value = 4453.7827
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13361.3481
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10911.3481
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 7292.5293
$ws.Range("I132").Value = 5838.3335
$ws.Range("J132").Value = 10782.6
$ws.Range("K132").Value = 17515.0005
$ws.Range("L132").Value = 32347.8
$ws.Range("M132").Value = -14985.0005
$ws.Range("N132").Value = -37407.8
$ws.Range("H136").Value = 4016.9814
$ws.Range("I136").Value = 3171.5
$ws.Range("J136").Value = 10780.833
$ws.Range("K136").Value = 9514.5
$ws.Range("L136").Value = 32342.499
$ws.Range("M136").Value = -6964.5
$ws.Range("N136").Value = -37442.499

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 196.125
$ws.Range("I4").Value = 305.77777
$ws.Range("J4").Value = 55.142857
$ws.Range("K4").Value = 305.77777
$ws.Range("L4").Value = 55.142857
$ws.Range("M4").Value = -190.77777
$ws.Range("N4").Value = -285.142857
$ws.Range("H20").Value = 2786
$ws.Range("I20").Value = 2506.75
$ws.Range("J20").Value = 3192.182
$ws.Range("K20").Value = 2506.75
$ws.Range("L20").Value = 3192.182
$ws.Range("M20").Value = -2259.75
$ws.Range("N20").Value = -3686.182
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H94").Value = 785.8
$ws.Range("I94").Value = 763.37933
$ws.Range("J94").Value = 894.1667
$ws.Range("K94").Value = 763.37933
$ws.Range("L94").Value = 894.1667
$ws.Range("M94").Value = -312.37933
$ws.Range("N94").Value = -1796.1667
$ws.Range("H134").Value = 7865
$ws.Range("I134").Value = 4326.4287
$ws.Range("J134").Value = 20250
$ws.Range("K134").Value = 12979.2861
$ws.Range("L134").Value = 60750
$ws.Range("M134").Value = -10444.2861
$ws.Range("N134").Value = -65820

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 454.41177
$ws.Range("I7").Value = 492.36
$ws.Range("K7").Value = 492.36
$ws.Range("M7").Value = -379.36
$ws.Range("H22").Value = 1354.875
$ws.Range("I22").Value = 686.7778
$ws.Range("J22").Value = 2213.8572
$ws.Range("K22").Value = 686.7778
$ws.Range("L22").Value = 2213.8572
$ws.Range("M22").Value = -336.7778
$ws.Range("N22").Value = -2913.8572
$ws.Range("H31").Value = 20002772
$ws.Range("I31").Value = 32259966
$ws.Range("J31").Value = 4194.3687
$ws.Range("K31").Value = 32259966
$ws.Range("L31").Value = 4194.3687
$ws.Range("M31").Value = -32259671
$ws.Range("N31").Value = -4784.3687
$ws.Range("H34").Value = 20002772
$ws.Range("I34").Value = 32259966
$ws.Range("J34").Value = 4194.3687
$ws.Range("K34").Value = 32259966
$ws.Range("L34").Value = 4194.3687
$ws.Range("M34").Value = -32259764
$ws.Range("N34").Value = -4598.3687
$ws.Range("H36").Value = 11024
$ws.Range("I36").Value = 11024
$ws.Range("K36").Value = 11024
$ws.Range("M36").Value = -10636
$ws.Range("H40").Value = 11024
$ws.Range("I40").Value = 11024
$ws.Range("K40").Value = 11024
$ws.Range("M40").Value = -10864
$ws.Range("H86").Value = 7043.2
$ws.Range("J86").Value = 6362.4
$ws.Range("L86").Value = 6362.4
$ws.Range("N86").Value = -8608.4
$ws.Range("H89").Value = 7043.2
$ws.Range("J89").Value = 6362.4
$ws.Range("L89").Value = 31812
$ws.Range("N89").Value = -43044
$ws.Range("H94").Value = 4881.077
$ws.Range("I94").Value = 8866.833000000001
$ws.Range("K94").Value = 8866.833000000001
$ws.Range("M94").Value = -8415.833000000001
$ws.Range("H132").Value = 79532.05
$ws.Range("I132").Value = 82397.13
$ws.Range("K132").Value = 247191.39
$ws.Range("M132").Value = -244661.39
$ws.Range("H134").Value = 7820.2354
$ws.Range("I134").Value = 6449.769
$ws.Range("K134").Value = 19349.307
$ws.Range("M134").Value = -16814.307
$ws.Range("H140").Value = 115000
$ws.Range("J140").Value = 115000
$ws.Range("L140").Value = 115000
$ws.Range("N140").Value = -125360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 681.25
$ws.Range("J23").Value = 681.25
$ws.Range("L23").Value = 2043.75
$ws.Range("N23").Value = -2513.75
$ws.Range("H39").Value = 5936.75
$ws.Range("J39").Value = 4872.5
$ws.Range("L39").Value = 14617.5
$ws.Range("N39").Value = -15205.5
$ws.Range("H55").Value = 1041.6666
$ws.Range("J55").Value = 1041.6666
$ws.Range("L55").Value = 3124.9998
$ws.Range("N55").Value = -3478.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9331.833000000001
$ws.Range("J70").Value = 9498
$ws.Range("L70").Value = 9498
$ws.Range("N70").Value = -10038
$ws.Range("H73").Value = 9331.833000000001
$ws.Range("J73").Value = 9498
$ws.Range("L73").Value = 9498
$ws.Range("N73").Value = -11370
$ws.Range("H97").Value = 1774.7
$ws.Range("I97").Value = 1815.75
$ws.Range("K97").Value = 1815.75
$ws.Range("M97").Value = -1319.75
$ws.Range("H122").Value = 26762.195
$ws.Range("I122").Value = 25315.643
$ws.Range("J122").Value = 29877.846
$ws.Range("K122").Value = 75946.929
$ws.Range("L122").Value = 89633.538
$ws.Range("M122").Value = -73496.929
$ws.Range("N122").Value = -94533.538
$ws.Range("H132").Value = 4986.162
$ws.Range("I132").Value = 4205.5293
$ws.Range("J132").Value = 13833.333
$ws.Range("K132").Value = 12616.5879
$ws.Range("L132").Value = 41499.999
$ws.Range("M132").Value = -10086.5879
$ws.Range("N132").Value = -46559.999
$ws.Range("H135").Value = 67915.75
$ws.Range("J135").Value = 67915.75
$ws.Range("L135").Value = 67915.75
$ws.Range("N135").Value = -78055.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1689.909
$ws.Range("I61").Value = 1776.6666
$ws.Range("J61").Value = 1299.5
$ws.Range("K61").Value = 1776.6666
$ws.Range("L61").Value = 1299.5
$ws.Range("M61").Value = -1574.6666
$ws.Range("N61").Value = -1703.5
$ws.Range("H82").Value = 2088.2222
$ws.Range("I82").Value = 1970.7142
$ws.Range("J82").Value = 2499.5
$ws.Range("K82").Value = 1970.7142
$ws.Range("L82").Value = 2499.5
$ws.Range("M82").Value = -1609.7142
$ws.Range("N82").Value = -3221.5
$ws.Range("H85").Value = 2088.2222
$ws.Range("I85").Value = 1970.7142
$ws.Range("J85").Value = 2499.5
$ws.Range("K85").Value = 1970.7142
$ws.Range("L85").Value = 2499.5
$ws.Range("M85").Value = -722.7141999999999
$ws.Range("N85").Value = -4995.5
$ws.Range("H93").Value = 1245.4117
$ws.Range("I93").Value = 1427.9286
$ws.Range("J93").Value = 393.66666
$ws.Range("K93").Value = 1427.9286
$ws.Range("L93").Value = 393.66666
$ws.Range("M93").Value = -179.9286
$ws.Range("N93").Value = -2889.66666
$ws.Range("H113").Value = 1689.909
$ws.Range("I113").Value = 1776.6666
$ws.Range("J113").Value = 1299.5
$ws.Range("K113").Value = 1776.6666
$ws.Range("L113").Value = 1299.5
$ws.Range("M113").Value = 393.3334
$ws.Range("N113").Value = -5639.5
$ws.Range("H122").Value = 8777
$ws.Range("I122").Value = 8777
$ws.Range("K122").Value = 26331
$ws.Range("M122").Value = -23881
$ws.Range("H136").Value = 3866.7646
$ws.Range("I136").Value = 2610
$ws.Range("K136").Value = 7830
$ws.Range("M136").Value = -5280

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3815.2917
$ws.Range("I107").Value = 4122.353
$ws.Range("J107").Value = 3069.5715
$ws.Range("K107").Value = 12367.059
$ws.Range("L107").Value = 9208.7145
$ws.Range("M107").Value = -10447.059
$ws.Range("N107").Value = -13048.7145
$ws.Range("H122").Value = 2438.9167
$ws.Range("I122").Value = 2226.3
$ws.Range("K122").Value = 6678.900000000001
$ws.Range("M122").Value = -4228.900000000001
$ws.Range("H132").Value = 4275.9385
$ws.Range("I132").Value = 3042.84
$ws.Range("K132").Value = 9128.52
$ws.Range("M132").Value = -6598.52
$ws.Range("H136").Value = 2728.8262
$ws.Range("I136").Value = 1066.5264
$ws.Range("J136").Value = 10624.75
$ws.Range("K136").Value = 3199.5792
$ws.Range("L136").Value = 31874.25
$ws.Range("M136").Value = -649.5792000000001
$ws.Range("N136").Value = -36974.25
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
